$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.012.93"
$ws.Range("E2").Value = "  -3.76%  "

$ws.Range("D3").Value = "1.650.69"
$ws.Range("E3").Value = "  -5.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.57"
$ws.Range("E5").Value = "  -5.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4832"
$ws.Range("E7").Value = "  -6.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2618"
$ws.Range("E8").Value = "  -5.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06014"
$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07193"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").Value = "1.655.33"
$ws.Range("E11").Value = "  -5.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.76"
$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6218"
$ws.Range("E13").Value = "  -4.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.565"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "73.02"
$ws.Range("E15").Value = "  -6.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9990"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "25.018.23"
$ws.Range("E18").Value = "  -3.80%  "

$ws.Range("E19").Value = "  -3.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006627"
$ws.Range("E20").Value = "  -2.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.556"
$ws.Range("E21").Value = "  +6.08%  "

$ws.Range("D22").Value = "1.857.25"
$ws.Range("E22").Value = "  -5.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.607"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("E25").Value = "  -2.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.96"
$ws.Range("E26").Value = "  -2.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("E27").Value = "  -7.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "103.02"
$ws.Range("E28").Value = "  -2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.681"
$ws.Range("E29").Value = "  -5.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.764"
$ws.Range("E30").Value = "  -4.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07896"
$ws.Range("E31").Value = "  -4.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.587"
$ws.Range("E32").Value = "  -2.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04554"
$ws.Range("E33").Value = "  -2.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9990"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.596"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9355"
$ws.Range("E36").Value = "  -6.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5799"
$ws.Range("E37").Value = "  -7.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.596"
$ws.Range("E38").Value = "  -5.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01564"
$ws.Range("E39").Value = "  -3.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8437"
$ws.Range("E40").Value = "  +10.39%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.827"
$ws.Range("E42").Value = "  -4.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.59"
$ws.Range("E43").Value = "  -1.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3731"
$ws.Range("E44").Value = "  -3.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.790"
$ws.Range("E45").Value = "  -4.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1150"
$ws.Range("E46").Value = "  +1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.157"
$ws.Range("E47").Value = "  -3.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05201"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.87"
$ws.Range("E49").Value = "  -2.95%  "

$ws.Range("E50").Value = "  -0.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.40"
$ws.Range("E51").Value = "  -9.48%  "
